$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra data rows (rows 4-9), keeping only header + rows 2 and 3
$ws.Rows("4:9").Delete() | Out-Null

# Insert a new column before "reward" (currently column G) for "lockdown_freq"
$ws.Columns("G:G").Insert() | Out-Null

# Insert a new column before "testing" (currently column I after the previous insert) for "test_freq"
$ws.Columns("I:I").Insert() | Out-Null

# --- Headers ---
$ws.Range("G1").Value = "lockdown_freq"
$ws.Range("I1").Value = "test_freq"

# --- Row 2 ---
$ws.Range("B2").Value = 11722.24006150288
$ws.Range("D2").Value = 8954554373.270102
$ws.Range("F2").Value = 3000
$ws.Range("G2").Value = 14
$ws.Range("H2").Value = -3831779429.893841
$ws.Range("I2").Value = 7
$ws.Range("K2").Value = 30000
$ws.Range("L2").Value = 1000000

# --- Row 3 ---
$ws.Range("B3").Value = 9833.931188451937
$ws.Range("D3").Value = 8956426690.199156
$ws.Range("F3").Value = 3000
$ws.Range("G3").Value = 14
$ws.Range("H3").Value = -1578617800.467741
$ws.Range("I3").Value = 7
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 1000000
